# Uppdaterat PPT med tillståndsmaskinstext
# Slide 23 ("Tillståndsmaskin: Equipped items (class Hero)"), content placeholder:
# replace the old "Motivering" paragraph with three new paragraphs describing
# the Equipment variables on Hero.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(23)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$para1 = "Equipment variablerna på Hero"
$para2text = "Vapen och Armor kan vara tilldelat eller null, skiljt från varandra, ger 4 olika tillstånd (nästkommande bild)."
$para3 = " 12 bågar då olika beteenden finns när man plockar upp ett item som är svagare eller starkare."

# Set the full text first, using carriage returns to create the three
# separate paragraphs that PowerPoint would produce for this edit.
$tr.Text = $para1 + "`r" + $para2text + "`r" + $para3

# Paragraph 2 needs to be split into separate runs (mirrors the run
# boundaries of the authored slide: "Armor" and "null" are their own runs).
$para2 = $tr.Paragraphs(2, 1)

$runArmor = $para2.Characters(11, 5)
$runArmor.Text = "Armor"

$runNull = $para2.Characters(42, 4)
$runNull.Text = "null"
